$p = $ppt.ActivePresentation
$s = $p.Slides.Item(15)
$shape = $s.Shapes.Item(2)
$tbl = $shape.Table

# Row 3 (Capstone Project 02), Weight column: 30 -> 20
$tbl.Cell(3, 2).Shape.TextFrame.TextRange.Text = "20"

# Row 4 (Case Study* (5X)), Component column: (5X -> (4X, keep the trailing ")" run intact
$tbl.Cell(4, 1).Shape.TextFrame.TextRange.Text = "Case Study* (4X"

# Row 4 (Case Study* (5X)), Weight column: 50 -> 60
$tbl.Cell(4, 2).Shape.TextFrame.TextRange.Text = "60"
